# Update "想去人数" (want-to-go count) values in column F across the
# three sheets that list event rows ("展览", "演出", "全部类型").
# "本地生活" is left untouched, as the diff does not modify it.

$wb = $excel.ActiveWorkbook

function Set-FValues {
    param(
        [string]$SheetName,
        [hashtable]$Updates
    )

    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $Updates.Keys) {
        $ws.Range("F$row").Value = $Updates[$row]
    }
}

# Sheet "展览" (sheet1)
Set-FValues "展览" @{
    4  = 5500
    5  = 5500
    6  = 217
    9  = 1240
    11 = 6
    12 = 810
    13 = 23
    14 = 6570
    17 = 125
    18 = 4858
    20 = 263
    21 = 4200
    23 = 4144
    25 = 216
    26 = 285
    27 = 266
    28 = 222
    33 = 61
    34 = 7506
    35 = 41
    36 = 1264
    37 = 623
    38 = 113
    39 = 980
    41 = 1495
    42 = 198
    43 = 833
    45 = 3679
    46 = 334
    47 = 15
    50 = 1038
}

# Sheet "演出" (sheet2)
Set-FValues "演出" @{
    10 = 9
    12 = 26
    13 = 151
    18 = 69
    21 = 853
}

# Sheet "全部类型" (sheet4)
Set-FValues "全部类型" @{
    6  = 5500
    7  = 5500
    8  = 217
    12 = 1240
    14 = 6
    16 = 810
    17 = 6570
    20 = 125
    21 = 4858
    23 = 263
    24 = 4200
    25 = 4144
    27 = 216
    28 = 285
    29 = 266
    30 = 222
    33 = 151
    34 = 7506
    35 = 41
    36 = 1264
    37 = 623
    38 = 113
    39 = 980
    41 = 1495
    42 = 198
    43 = 833
    45 = 3679
    46 = 334
    49 = 1038
}
